$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$html = @"

    <!DOCTYPE html>
    <html>
    <head>
        <style>
            table {
                font-family: Arial, sans-serif;
                border-collapse: collapse;
                width: 100%;
            }
            th, td {
                border: 1px solid #dddddd;
                text-align: left;
                padding: 8px;
            }
            th {
                background-color: #f2f2f2;
            }
        </style>
    </head>
    <body>
    <h2>Error Log</h2>
    <table>
        <tr>
            <th>Timestamp</th>
            <td>2024-11-04 10:44:36</td>
        </tr>
        <tr>
            <th>Error Level</th>
            <td>High</td>
        </tr>
        <tr>
            <th>Location</th>
            <td>tasks.py, line 22</td>
        </tr>
        <tr>
            <th>Error Message</th>
            <td>Testing!!!</td>
        </tr>
    </table>
    </body>
    </html>
    
"@

$ws.Range("A30").Value = "2024-11-04 10:44:40"
$ws.Range("B30").Value = "Success"
$ws.Range("C30").Value = $html

# The multi-line HTML content otherwise triggers an auto-calculated explicit
# row height (customHeight="1"); AutoFit clears that flag so the row matches
# the target sheet's un-pinned height, same as row 29's existing pattern.
$ws.Rows(30).AutoFit()
